# Updates cryptos list values (price/volume) to reflect the latest scrape.
# Cells in column D whose new text is a plain numeric string are entered
# with a leading apostrophe so Excel stores them as text (matching the
# workbook's existing inlineStr convention) instead of auto-converting
# them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.988.40"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.515.78"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").Value = "'532.48"
$ws.Range("E5").Value = "  -2.09%  "
$ws.Range("D6").Value = "'138.56"
$ws.Range("E6").Value = "  -4.44%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -2.24%  "
$ws.Range("D9").Value = "2.518.26"
$ws.Range("E9").Value = "  -1.38%  "
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("D13").Value = "'0.355"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "2.957.90"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").Value = "'23.15"
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").Value = "58.928.47"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "'0.0000140"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "2.515.05"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").Value = "'11.04"
$ws.Range("E19").Value = "  -2.17%  "
$ws.Range("D20").Value = "'4.29"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "'323.13"
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'5.81"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").Value = "'62.12"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  -3.49%  "
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  -2.83%  "
$ws.Range("D29").Value = "'6.69"
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("D30").Value = "0.0₃0770"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("D31").Value = "'1.79"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("D32").Value = "'164.39"
$ws.Range("E32").Value = "  +4.95%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.45"
$ws.Range("E34").Value = "  -3.49%  "
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.12"
$ws.Range("E35").Value = "  -9.10%  "
$ws.Range("D36").Value = "'18.48"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").Value = "'4.24"
$ws.Range("E37").Value = "  -3.80%  "
$ws.Range("E38").Value = "  -3.39%  "
$ws.Range("D39").Value = "'36.76"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").Value = "'3.66"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").Value = "'0.807"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("D42").Value = "'5.22"
$ws.Range("E42").Value = "  -8.43%  "
$ws.Range("D43").Value = "'278.82"
$ws.Range("E43").Value = "  -7.31%  "
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").Value = "'10.88"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'122.14"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("D49").Value = "'18.43"
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "'0.0511"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").Value = "'0.0224"
$ws.Range("E51").Value = "  -2.29%  "
